# Applies the Fri Oct 18 16:13:18 UTC 2024 "Updated cryptos list" GitHub
# Actions refresh: per-coin Price (D) / Volume(1h) (E) updates, plus the OKB<->USDe
# and Optimism<->InjectiveProtocol row swaps (rows 46/47 and 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.550.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.79%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''2.637.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.55%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '''  +0.03%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''599.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.17%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''154.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +1.49%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '''  +0.05%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.544'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -1.39%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''2.641.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.69%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''0.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +11.93%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = '''  -0.49%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = '''  +0.94%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''0.348'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.25%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''27.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.25%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''0.0000188'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +5.30%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''3.124.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.18%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''68.494.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +2.06%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''2.632.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.35%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''11.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.96%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''366.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.31%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''7.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +1.00%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''4.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.01%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''4.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.09%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''2.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.66%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''73.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.43%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = '''  -0.06%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''9.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.81%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''2.790.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +1.10%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''0.0000104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +3.51%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = '''  +0.01%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''577.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -1.08%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = '''7.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +3.01%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''1.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.23%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''1.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.14%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = '''  +0.02%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = '''0.127'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +3.72%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = '''  +1.79%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''160.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +2.66%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''19.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.16%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = '''  +2.16%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''0.366'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.15%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''5.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +1.96%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''2.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +2.64%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''17.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +5.40%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''0.0₆0320'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +8.22%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = '''USDe'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.04%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = '''OKB'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''40.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.97%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''155.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.28%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''3.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -0.33%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = '''InjectiveProtocol'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''21.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +0.16%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = '''Optimism'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''1.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -0.07%  '
$ws.Range("E51").Style = "Normal"

